$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "总计" (Total) summary sheet
$ws2 = $wb.Worksheets.Item(2)   # "2022-Q2" sheet -> becomes "2022-Q3"

# Duplicate ws2 right after itself. The duplicate keeps the original
# "2022-Q2" data / margins / styles untouched and becomes the new
# "2022-Q2" sheet (sheetId 3), while ws2 itself (sheetId 2) is turned
# into the new "2022-Q3" sheet.
$ws2.Copy($null, $ws2)
$wsQ2 = $wb.Worksheets.Item(3)

$ws2.Name = "2022-Q3"
$wsQ2.Name = "2022-Q2"

# The freshly-populated "2022-Q3" sheet uses the same page margins as
# the "总计" sheet (1in top/bottom, 0.75in left/right, 0.5in header/footer)
# rather than the old sheet's margins, so fix those up explicitly.
# PageSetup margins are expressed in points (72pt = 1in).
$ws2.PageSetup.LeftMargin = 0.75 * 72
$ws2.PageSetup.RightMargin = 0.75 * 72
$ws2.PageSetup.TopMargin = 1 * 72
$ws2.PageSetup.BottomMargin = 1 * 72
$ws2.PageSetup.HeaderMargin = 0.5 * 72
$ws2.PageSetup.FooterMargin = 0.5 * 72

# Header row (B1:H1) and the A column on "2022-Q3" use the same style as
# the "总计" sheet's header / A column, so copy that formatting over.
$ws1.Range("B1").Copy()
$ws2.Range("B1:H1").PasteSpecial(-4122)
$ws1.Range("A2").Copy()
$ws2.Range("A2:A3").PasteSpecial(-4122)

function Set-TextValue($range, $value) {
    # Force the value to be stored as text (matching the source data,
    # which keeps values such as "15.80" / "202801" as text, not
    # numbers), then strip the number-format styling that this requires
    # so the cell is left without any extra style applied.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

# Row 2 of "2022-Q3": fund 202801
$ws2.Range("A2").Value = 0
Set-TextValue $ws2.Range("B2") "202801"
Set-TextValue $ws2.Range("C2") "南方全球精选配置（QDII-FOF）"
Set-TextValue $ws2.Range("D2") "15.80"
Set-TextValue $ws2.Range("E2") "29.52"
Set-TextValue $ws2.Range("F2") "1.04"
Set-TextValue $ws2.Range("G2") "0.1643"
$ws2.Range("H2").Value = 10

# Row 3 of "2022-Q3": fund 519602
$ws2.Range("A3").Value = 1
Set-TextValue $ws2.Range("B3") "519602"
Set-TextValue $ws2.Range("C3") "海富通大中华精选混合（QDII）"
Set-TextValue $ws2.Range("D3") "0.10"
Set-TextValue $ws2.Range("E3") "87.37"
Set-TextValue $ws2.Range("F3") "3.33"
Set-TextValue $ws2.Range("G3") "0.0033"
$ws2.Range("H3").Value = 10

# Update the "总计" (Total) summary sheet: row 2 now reports the new
# "2022-Q3" totals, and a new row 3 is added reporting the original
# "2022-Q2" totals (previously on row 2).
$ws1.Range("A2").Copy()
$ws1.Range("A3").PasteSpecial(-4122)
$ws1.Range("A3").Value = 1

$ws1.Range("B2").Value = "2022-Q3"
$ws1.Range("C2").Value = 2
$ws1.Range("D2").Value = 0.17

$ws1.Range("B3").Value = "2022-Q2"
$ws1.Range("C3").Value = 1
$ws1.Range("D3").Value = 0.07000000000000001

# Keep "总计" as the active sheet/tab, matching the source workbook.
$ws1.Activate()
